$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -10.85509999999999
$ws.Range("D3").Value = -7.206699999999993
$ws.Range("B12").Value = 4.929499999999996
$ws.Range("C14").Value = -13.3753
$ws.Range("C26").Value = -11.9852
$ws.Range("D30").Value = -7.489999999999998
$ws.Range("C31").Value = -13.6568
$ws.Range("B32").Value = 6.418599999999999
$ws.Range("C35").Value = -12.66540000000002
$ws.Range("B36").Value = 8.983900000000004
$ws.Range("C37").Value = -12.9902
$ws.Range("B38").Value = 5.815
$ws.Range("D44").Value = -7.315400000000003
$ws.Range("C45").Value = -13.4973
$ws.Range("B46").Value = 6.515300000000002
$ws.Range("B54").Value = 4.702300000000002
$ws.Range("B55").Value = 5.355899999999998
$ws.Range("C57").Value = -14.23859999999999
$ws.Range("D58").Value = -7.944899999999993
$ws.Range("B67").Value = 5.116399999999995
$ws.Range("B69").Value = 4.948699999999995
$ws.Range("B72").Value = 5.3531
$ws.Range("D84").Value = -8.752600000000003
$ws.Range("D89").Value = -6.096099999999998
$ws.Range("B91").Value = 5.0285
$ws.Range("D91").Value = -6.134099999999997
$ws.Range("D92").Value = -6.004400000000001
$ws.Range("B99").Value = 4.313299999999999
$ws.Range("C100").Value = -12.99019999999999
$ws.Range("C102").Value = -13.12040000000001
$ws.Range("D102").Value = -7.715499999999996
